$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('F130').Value = 'https://www.perekrestok.ru/cat/391/p/kofe-nescafe-gold-naturalnyj-rastvorimyj-s-dobavleniem-molotogo-95g-3659496'
$ws.Range('F131').Value = 'https://www.perekrestok.ru/cat/391/p/kofe-nescafe-gold-naturalnyj-rastvorimyj-s-dobavleniem-molotogo-190g-3928172'
$ws.Range('F132').Value = 'https://www.perekrestok.ru/cat/391/p/kofe-nescafe-gold-barista-naturalnyj-rastvorimyj-s-dobavleniem-molotogo-85g-3419749'
$ws.Range('F134').Value = 'https://www.perekrestok.ru/cat/391/p/kofe-nescafe-classic-rastvorimyj-s-dobavleniem-molotogo-95g-4035358'
$ws.Range('F135').Value = 'https://www.perekrestok.ru/cat/391/p/kofe-nescafe-gold-crema-naturalnyj-rastvorimyj-poroskoobraznyj-95g-3626842'
$ws.Range('F136').Value = 'https://www.perekrestok.ru/cat/403/p/marmelad-bon-pari-kislye-cervacki-zevatelnyj-frukty-75g-43332'
$ws.Range('F137').Value = 'https://www.perekrestok.ru/cat/403/p/marmelad-bon-pari-masiny-zevatelnyj-frukty-75g-2127571'
$ws.Range('F138').Value = 'https://www.perekrestok.ru/cat/446/p/vlaznyj-korm-felix-dla-kosek-s-kuricej-85g-3331910'
$ws.Range('F139').Value = 'https://www.perekrestok.ru/cat/446/p/vlaznyj-korm-felix-dla-kosek-s-lososem-85g-3413297'
$ws.Range('F140').Value = 'https://www.perekrestok.ru/cat/446/p/vlaznyj-korm-felix-dla-kosek-s-govadinoj-85g-3356712'
$ws.Range('F141').Value = 'https://www.perekrestok.ru/cat/446/p/vlaznyj-korm-felix-dla-kosek-s-kuricej-48g-3931384'
$ws.Range('F142').Value = 'https://www.perekrestok.ru/cat/446/p/vlaznyj-korm-purina-one-dla-kosek-s-kuricej-i-zelenoj-fasolu-75g-4002292'
$ws.Range('F143').Value = 'https://www.perekrestok.ru/cat/446/p/korm-suhoj-felix-veselaa-ohota-so-vkusom-krevetki-i-ryby-dla-kosek-50g-4090222'
$ws.Range('F144').Value = 'https://www.perekrestok.ru/cat/393/p/maslo-podsolnecnoe-sloboda-rafinirovannoe-dezodorirovannoe-vyssij-sort-1l-21676'
$ws.Range('F146').Value = 'https://www.perekrestok.ru/cat/393/p/maslo-podsolnecnoe-sloboda-rafinirovannoe-1-8l-2147548'
$ws.Range('F148').Value = 'https://www.perekrestok.ru/cat/394/p/majonez-sloboda-olivkovyj-67-400ml-2108609'
$ws.Range('F149').Value = 'https://www.perekrestok.ru/cat/394/p/majonez-sloboda-olivkovyj-67-800ml-3163111'
$ws.Range('F150').Value = 'https://www.perekrestok.ru/cat/394/p/majonez-sloboda-provansal-olivkovyj-67-230ml-3389322'
$ws.Range('F151').Value = 'https://www.perekrestok.ru/cat/394/p/majonez-sloboda-provansal-67-400ml-3674577'
$ws.Range('F152').Value = 'https://www.perekrestok.ru/cat/394/p/majonez-sloboda-provansal-67-800ml-3193489'
$ws.Range('F153').Value = 'https://www.perekrestok.ru/cat/394/p/majonez-sloboda-provansal-67-230ml-3389323'
$ws.Range('F154').Value = 'https://www.perekrestok.ru/cat/394/p/sous-majoneznyj-sloboda-postnyj-56-200ml-3311552'
$ws.Range('F156').Value = 'https://www.perekrestok.ru/cat/394/p/sous-sloboda-gribnoj-60-220ml-3483228'
$ws.Range('F157').Value = 'https://www.perekrestok.ru/cat/394/p/ketcup-sloboda-tomatnyj-350g-3382584'
$ws.Range('F158').Value = 'https://www.perekrestok.ru/cat/372/p/biojogurt-sloboda-s-granola-i-izumom-obogasennyj-laktobakteriami-5-130g-4204914'
$ws.Range('F160').Value = 'https://www.perekrestok.ru/cat/372/p/jogurt-liberty-yogurt-s-dynej-i-dikim-risom-1-5-270g-4212430'
$ws.Range('F161').Value = 'https://www.perekrestok.ru/cat/372/p/jogurt-liberty-yogurt-s-arbuzom-matoj-i-bambukom-1-5-270g-4212429'
$ws.Range('F162').Value = 'https://www.perekrestok.ru/cat/390/p/lapsa-rollton-aicnaa-klassiceskaa-400g-45105'
$ws.Range('F163').Value = 'https://www.perekrestok.ru/cat/465/p/pure-rollton-goracaa-porcia-kartofelnoe-s-kuricej-55g-3296248'
$ws.Range('F164').Value = 'https://www.perekrestok.ru/cat/465/p/lapsa-aicnaa-rollton-kurinaa-po-domasnemu-bystrogo-prigotovlenia-90g-3453979'
$ws.Range('F165').Value = 'https://www.perekrestok.ru/cat/465/p/lapsa-rollton-s-govadinoj-po-domasnemu-90g-3453968'
$ws.Range('F166').Value = 'https://www.perekrestok.ru/cat/465/p/bulon-rollton-kurinyj-domasnij-100g-76002'
$ws.Range('F168').Value = 'https://www.perekrestok.ru/cat/465/p/pure-kartofelnoe-big-bon-po-domasnemu-s-kusockami-tusenoj-govadiny-110g-4144382'
$ws.Range('G127').Value = 'https://lavka.yandex.ru/213/good/b60483e97ee6496b9dee0bf8d6e263c0000100010001?searchQuery=%D0%BD%D0%B5%D1%81%D0%BA%D0%B0%D1%84%D0%B5&searchPosition=-1'
$ws.Range('G128').Value = 'https://lavka.yandex.ru/213/good/04c715bb8ba0487e8e5e70034e872932000100010001?searchQuery=%D0%BD%D0%B5%D1%81%D0%BA%D0%B0%D1%84%D0%B5&searchPosition=-1'
$ws.Range('G130').Value = 'https://lavka.yandex.ru/213/good/a40d2103fddd4583a7ef7967aea8fa57000300010000?searchQuery=%D1%84%D0%B5%D0%BB%D0%B8%D0%BA%D1%81&searchPosition=-1'
$ws.Range('G131').Value = 'https://lavka.yandex.ru/213/good/f54177c70ddb417faec89e35f6de766c000200010001?searchQuery=one&searchPosition=-1'
$ws.Range('G132').Value = 'https://lavka.yandex.ru/213/good/774d01e6c4a143a3b188ac2b658f65aa000200010001?searchQuery=%D0%B0%D0%BB%D1%8C%D1%82%D0%B5%D1%80%D0%BE&searchPosition=-1'
$ws.Range('G134').Value = 'https://lavka.yandex.ru/213/good/ff9cdb35c790411faf9dce88544afd73000100010001?searchQuery=%D0%B0%D0%BB%D1%8C%D1%82%D0%B5%D1%80%D0%BE&searchPosition=-1'
$ws.Range('G133').Value = 'https://lavka.yandex.ru/213/good/2fec605ca88f449abf6f88b0a3879d6c000300010000?searchQuery=%D1%81%D0%BB%D0%BE%D0%B1%D0%BE%D0%B4%D0%B0&searchPosition=-1'
$ws.Range('G135').Value = 'https://lavka.yandex.ru/213/good/1c68511883fc4325b46794f87a3df06a000300010000?searchQuery=%D1%81%D0%BB%D0%BE%D0%B1%D0%BE%D0%B4%D0%B0&searchPosition=-1'
$ws.Range('G136').Value = 'https://lavka.yandex.ru/213/good/a3e1fcfb1c0c4fb1b24f40b49a821659000200010001?searchQuery=%D1%81%D0%BB%D0%BE%D0%B1%D0%BE%D0%B4%D0%B0&searchPosition=-1'
$ws.Range('G137').Value = 'https://lavka.yandex.ru/213/good/309ac27abceb436b984f44fbe7e4989b000300010000?searchQuery=%D1%81%D0%BB%D0%BE%D0%B1%D0%BE%D0%B4%D0%B0&searchPosition=-1'
$ws.Range('G138').Value = 'https://lavka.yandex.ru/213/good/7cf94231b92a4aa0b0cd0d7004ec9d7f000200010000?searchQuery=%D1%81%D0%BB%D0%BE%D0%B1%D0%BE%D0%B4%D0%B0&searchPosition=-1'
$ws.Range('G139').Value = 'https://lavka.yandex.ru/213/good/d1b5202b49b0402184f181c04b271b65000100010001?searchQuery=%D1%81%D0%BB%D0%BE%D0%B1%D0%BE%D0%B4%D0%B0&searchPosition=-1'
$ws.Range('G140').Value = 'https://lavka.yandex.ru/213/good/0b8ff53d6a3546de9c76aae8076aafca000300010000?searchQuery=%D1%81%D0%BB%D0%BE%D0%B1%D0%BE%D0%B4%D0%B0&searchPosition=-1'
$ws.Range('G141').Value = 'https://lavka.yandex.ru/213/good/8b681e16f31a4ea7b21b03979b9e6501000100010001?searchQuery=%D1%81%D0%BB%D0%BE%D0%B1%D0%BE%D0%B4%D0%B0&searchPosition=-1'
$ws.Range('G142').Value = 'https://lavka.yandex.ru/213/good/104adeac520045c88d806c402b969498000200010000?searchQuery=%D1%80%D0%BE%D0%BB%D0%BB%D1%82%D0%BE%D0%BD%20&searchPosition=-1'
$ws.Range('G143').Value = 'https://lavka.yandex.ru/213/good/58272528d78e4a8f944e258bd5bb4fea000100010001?searchQuery=%D1%80%D0%BE%D0%BB%D0%BB%D1%82%D0%BE%D0%BD%20&searchPosition=-1'
$ws.Range('G144').Value = 'https://lavka.yandex.ru/213/good/c43a2c9ab06343f3b36670f1c013d8f6000200010000?searchQuery=%D0%B1%D0%B8%D0%B3%20%D1%8E%D0%B1%D0%BE%D0%BD&searchPosition=-1'
$ws.Range('J140').Value = 'https://www.utkonos.ru/item/3347661/kofe-nescafe-gold-rastvorimyj-sublimirovannyj-95-g'
$ws.Range('J141').Value = 'https://www.utkonos.ru/item/3347653/kofe-nescafe-gold-rastvorimyj-sublimirovannyj-190-g'
$ws.Range('J142').Value = 'https://www.utkonos.ru/item/3351671/kofe-nescafe-gold-barista-rastvorimyj-sublimirovannyj-85-g'
$ws.Range('J143').Value = 'https://www.utkonos.ru/item/3399235/kofe-nescafe-classic-rastvorimyj-poroshkoobraznyj-95-g'
$ws.Range('J144').Value = 'https://www.utkonos.ru/item/3125340/marmelad-bon-pari-kislye-chervjachki-zhevatelnyj-75-g'
$ws.Range('J145').Value = 'https://www.utkonos.ru/item/3162862/marmelad-bon-pari-zhevatelnyj-mashinki-so-vkusom-fruktov-75-g'
$ws.Range('J146').Value = 'https://www.utkonos.ru/item/3142734/korm-vlazhnyj-felix-v-zhele-s-kuricej-dlja-koshek-85-g'
$ws.Range('J147').Value = 'https://www.utkonos.ru/item/3290126/korm-vlazhnyj-felix-appetitnye-kusochki-dvojnaja-vkusnjatina-s-lososem-i-forelju-dlja-vzroslykh-koshek-85-g'
$ws.Range('J148').Value = 'https://www.utkonos.ru/item/3383786/korm-sukhoj-felix-dvojnaja-vkusnjatina-s-mjasom-dlja-vzroslykh-koshek-300-g'
$ws.Range('J149').Value = 'https://www.utkonos.ru/item/3356811/maslo-altero-podsolnechnoe-golden-s-dobavleniem-olivkovogo-dezodorirovannoe-rafinirovannoe-810-ml'
$ws.Range('J151').Value = 'https://www.utkonos.ru/item/3356844/maslo-altero-kukuruznoe-beauty-dezodorirovannoe-rafinirovannoe-810-ml'
$ws.Range('J150').Value = 'https://www.utkonos.ru/item/3356913/maslo-sloboda-podsolnechnoe-dezodorirovannoe-rafinirovannoe-1-l'
$ws.Range('J152').Value = 'https://www.utkonos.ru/item/3347975/maslo-sloboda-podsolnechnoe-aromatnoe-nerafinirovannoe-1-l'
$ws.Range('J153').Value = 'https://www.utkonos.ru/item/3402884/maslo-altero-podsolnechnoe-almond-rafinirovannoe-s-dobavleniem-mindalnogo-masla-810-ml'
$ws.Range('J154').Value = 'https://www.utkonos.ru/item/3052620/majonez-sloboda-olivkovyj-67-400-ml'
$ws.Range('J155').Value = 'https://www.utkonos.ru/item/3056060/majonez-sloboda-olivkovyj-67-800-ml'
$ws.Range('J156').Value = 'https://www.utkonos.ru/item/3050852/majonez-sloboda-olivkovyj-67-230-ml'
$ws.Range('J157').Value = 'https://www.utkonos.ru/item/3352751/majonez-sloboda-provansal-67-400-ml'
$ws.Range('J158').Value = 'https://www.utkonos.ru/item/3056059/majonez-sloboda-provansal-67-800-ml'
$ws.Range('J159').Value = 'https://www.utkonos.ru/item/3052658/majonez-sloboda-provansal-67-230-ml'
$ws.Range('J160').Value = 'https://www.utkonos.ru/item/3382065/majonez-sloboda-provansal-omega-3-6-9-67-400-ml'
$ws.Range('J161').Value = 'https://www.utkonos.ru/item/3428842/majoneznyj-sous-sloboda-postnyj-56-200-ml'
$ws.Range('J162').Value = 'https://www.utkonos.ru/item/3382073/sous-sloboda-cezar-220-ml'
$ws.Range('J163').Value = 'https://www.utkonos.ru/item/3382071/sous-sloboda-gribnoj-220-ml'
$ws.Range('J164').Value = 'https://www.utkonos.ru/item/3382063/ketchup-sloboda-tomatnyj-350-g'
$ws.Range('J165').Value = 'https://www.utkonos.ru/item/3382064/ketchup-sloboda-shashlychnyj-350-g'
$ws.Range('J166').Value = 'https://www.utkonos.ru/item/3436948/ketchup-sloboda-gurme-350-g'
$ws.Range('J167').Value = 'https://www.utkonos.ru/item/3411523/biojogurt-sloboda-semejnyj-zavtrak-granola-greckij-orekh-izjum-5-170-g'
$ws.Range('J168').Value = 'https://www.utkonos.ru/item/3411510/biojogurt-sloboda-pitevoj-vishnja-2-260-g'
$ws.Range('J169').Value = 'https://www.utkonos.ru/item/3411513/biojogurt-sloboda-pitevoj-persik-2-260-g'
$ws.Range('J170').Value = 'https://www.utkonos.ru/item/3411521/biojogurt-sloboda-semejnyj-zavtrak-mjusli-jabloko-mandarin-orekhi-5-170-g'
$ws.Range('J171').Value = 'https://www.utkonos.ru/item/3426857/biojogurt-sloboda-termostatnyj-mindalnyj-2-150-g'
$ws.Range('J172').Value = 'https://www.utkonos.ru/item/3103076/lapsha-rollton-jaichnaja-klassicheskaja-400-g'
$ws.Range('J173').Value = 'https://www.utkonos.ru/item/3391124/pjure-rollton-kartofelnoe-s-kuricej-55-g'
$ws.Range('J174').Value = 'https://www.utkonos.ru/item/3261763/lapsha-rollton-po-domashnemu-kurinaja-bystrogo-prigotovlenija-90-g'
$ws.Range('J175').Value = 'https://www.utkonos.ru/item/3416048/lapsha-rollton-s-kuricej-bystrogo-prigotovlenija-70-g'
$ws.Range('J176').Value = 'https://www.utkonos.ru/item/3261761/lapsha-rollton-po-domashnemu-s-govjadinoj-bystrogo-prigotovlenija-90-g'
$ws.Range('J177').Value = 'https://www.utkonos.ru/item/3416047/lapsha-big-bon-max-s-kuricej-bystrogo-prigotovlenija-95-g'
$ws.Range('J178').Value = 'https://www.utkonos.ru/item/3416049/lapsha-big-bon-max-s-govjadinoj-bystrogo-prigotovlenija-95-g'
$ws.Range('J179').Value = 'https://www.utkonos.ru/item/3391129/lapsha-big-bon-gotovyj-obed-s-kuricej-bystrogo-prigotovlenija-110-g'
